# Repull data, push all data, mean calculation
# Update the dSF (column F) values for the affected rows to match
# the freshly repulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = -2
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("F19").Value = -1
$ws.Range("F23").Value = 5
$ws.Range("F29").Value = 1
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = -5
$ws.Range("F35").Value = -1
$ws.Range("F41").Value = 5
$ws.Range("F42").Value = -4
